$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: CV Sent -> 1st Interview, date 45987 -> 45996
$ws.Range("E6").Value = "1st Interview"
$ws.Range("F6").Value = 45996

# Row 7: CV Sent -> 1st Interview, date 45987 -> 45993
$ws.Range("E7").Value = "1st Interview"
$ws.Range("F7").Value = 45993

# Row 8: CV Sent -> 1st Interview, date 45989 -> 45994
$ws.Range("E8").Value = "1st Interview"
$ws.Range("F8").Value = 45994

# Row 11: 2nd Interview -> 3rd Interview, date 45980 -> 45992
$ws.Range("E11").Value = "3rd Interview"
$ws.Range("F11").Value = 45992
